$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data (A1:B6) and rebuild A1:G4 with the new "book inventory" layout.
$ws.Range("A1:B6").Clear()

# Header row
$ws.Range("A1").Value = "Tên sách"
$ws.Range("B1").Value = "Mã nhà xuất bản"
$ws.Range("C1").Value = "mã tác giả"
$ws.Range("D1").Value = "mã thể loại"
$ws.Range("E1").Value = "số lượng tồn"
$ws.Range("F1").Value = "năm xuất bản"
$ws.Range("G1").Value = "đơn giá"

# Row 2
$ws.Range("A2").Value = "TestExcel1"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 123
$ws.Range("F2").Value = 2012
$ws.Range("G2").Value = 1000

# Row 3
$ws.Range("A3").Value = "TestExcel2"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 123
$ws.Range("F3").Value = 2012
$ws.Range("G3").Value = 1000

# Row 4
$ws.Range("A4").Value = "TestExcel3"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 123
$ws.Range("F4").Value = 2012
$ws.Range("G4").Value = 1000

# Column widths to match the authored layout (closest achievable given the
# host's character->pixel->character rounding of ColumnWidth)
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(3).ColumnWidth = 15.166666666666666
$ws.Columns.Item(4).ColumnWidth = 13.666666666666666
$ws.Columns.Item(5).ColumnWidth = 12.5
$ws.Columns.Item(6).ColumnWidth = 12.166666666666666
$ws.Columns.Item(7).ColumnWidth = 15.5

# Selection as left in the authored file
$ws.Range("I12").Select() | Out-Null
